$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy styles for the newly appended rows (207-215) from row 206 ---
$ws.Range("A206").Copy() | Out-Null
$ws.Range("A207:A215").PasteSpecial(-4122) | Out-Null
$ws.Range("E206").Copy() | Out-Null
$ws.Range("E207:E215").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 203
$ws.Range("A203").Value = 201
$ws.Range("B203").Value = 6774880
$ws.Range("C203").Value = 'Poland Ekstraklasa'
$ws.Range("D203").Value = 'Poland Ekstraklasa'
$ws.Range("E203").Value = 45354.35416666666
$ws.Range("F203").Value = 'LKS Lodz'
$ws.Range("G203").Value = 'MKS Puszcza Niepolomice'
$ws.Range("H203").Value = 3
$ws.Range("I203").Value = 2
$ws.Range("J203").Value = 'H'
$ws.Range("K203").Value = 2.45
$ws.Range("L203").Value = 3.3
$ws.Range("M203").Value = 2.55
$ws.Range("N203").Value = 2.25
$ws.Range("O203").Value = 3.3
$ws.Range("P203").Value = 2.875
$ws.Range("Q203").Value = -0.25
$ws.Range("R203").Value = 2.025
$ws.Range("S203").Value = 1.825
$ws.Range("T203").Value = 2.25
$ws.Range("U203").Value = 1.8
$ws.Range("V203").Value = 2.05
$ws.Range("W203").Value = 1.25
$ws.Range("X203").Value = -1
$ws.Range("Y203").Value = -1
$ws.Range("Z203").Value = 1.025
$ws.Range("AA203").Value = -1
$ws.Range("AB203").Value = 0.8
$ws.Range("AC203").Value = -1

# Row 204
$ws.Range("A204").Value = 202
$ws.Range("B204").Value = 6775556
$ws.Range("C204").Value = 'Poland Ekstraklasa'
$ws.Range("D204").Value = 'Poland Ekstraklasa'
$ws.Range("E204").Value = 45354.45833333334
$ws.Range("F204").Value = 'Zaglebie Lubin'
$ws.Range("G204").Value = 'Korona Kielce'
$ws.Range("H204").Value = 1
$ws.Range("I204").Value = 0
$ws.Range("J204").Value = 'H'
$ws.Range("K204").Value = 2.05
$ws.Range("L204").Value = 3.2
$ws.Range("M204").Value = 3.2
$ws.Range("N204").Value = 1.95
$ws.Range("O204").Value = 3.2
$ws.Range("P204").Value = 3.4
$ws.Range("Q204").Value = -0.5
$ws.Range("R204").Value = 2.05
$ws.Range("S204").Value = 1.8
$ws.Range("T204").Value = 2.5
$ws.Range("U204").Value = 2.05
$ws.Range("V204").Value = 1.8
$ws.Range("W204").Value = 0.95
$ws.Range("X204").Value = -1
$ws.Range("Y204").Value = -1
$ws.Range("Z204").Value = 1.05
$ws.Range("AA204").Value = -1
$ws.Range("AB204").Value = -1
$ws.Range("AC204").Value = 0.8

# Row 205
$ws.Range("A205").Value = 203
$ws.Range("B205").Value = 6775559
$ws.Range("C205").Value = 'Poland Ekstraklasa'
$ws.Range("D205").Value = 'Poland Ekstraklasa'
$ws.Range("E205").Value = 45354.5625
$ws.Range("F205").Value = 'Rakow Czestochowa'
$ws.Range("G205").Value = 'Lech Poznan'
$ws.Range("H205").Value = 4
$ws.Range("I205").Value = 0
$ws.Range("J205").Value = 'H'
$ws.Range("K205").Value = 1.95
$ws.Range("L205").Value = 3.25
$ws.Range("M205").Value = 3.8
$ws.Range("N205").Value = 1.95
$ws.Range("O205").Value = 3.25
$ws.Range("P205").Value = 3.75
$ws.Range("Q205").Value = -0.5
$ws.Range("R205").Value = 2.025
$ws.Range("S205").Value = 1.825
$ws.Range("T205").Value = 2.25
$ws.Range("U205").Value = 2
$ws.Range("V205").Value = 1.85
$ws.Range("W205").Value = 0.95
$ws.Range("X205").Value = -1
$ws.Range("Y205").Value = -1
$ws.Range("Z205").Value = 1.025
$ws.Range("AA205").Value = -1
$ws.Range("AB205").Value = 1
$ws.Range("AC205").Value = -1

# Row 206
$ws.Range("A206").Value = 204
$ws.Range("B206").Value = 6775558
$ws.Range("C206").Value = 'Poland Ekstraklasa'
$ws.Range("D206").Value = 'Poland Ekstraklasa'
$ws.Range("E206").Value = 45355.625
$ws.Range("F206").Value = 'Radomiak Radom'
$ws.Range("G206").Value = 'Stal Mielec'
$ws.Range("H206").Value = 2
$ws.Range("I206").Value = 1
$ws.Range("J206").Value = 'H'
$ws.Range("K206").Value = 2.05
$ws.Range("L206").Value = 3.2
$ws.Range("M206").Value = 3.5
$ws.Range("N206").Value = 1.833
$ws.Range("O206").Value = 3.2
$ws.Range("P206").Value = 4.5
$ws.Range("Q206").Value = -0.5
$ws.Range("R206").Value = 1.85
$ws.Range("S206").Value = 2
$ws.Range("T206").Value = 2.25
$ws.Range("U206").Value = 2.05
$ws.Range("V206").Value = 1.8
$ws.Range("W206").Value = 0.833
$ws.Range("X206").Value = -1
$ws.Range("Y206").Value = -1
$ws.Range("Z206").Value = 0.8500000000000001
$ws.Range("AA206").Value = -1
$ws.Range("AB206").Value = 1.05
$ws.Range("AC206").Value = -1

# Row 207
$ws.Range("A207").Value = 205
$ws.Range("B207").Value = 6775564
$ws.Range("C207").Value = 'Poland Ekstraklasa'
$ws.Range("D207").Value = 'Poland Ekstraklasa'
$ws.Range("E207").Value = 45359.58333333334
$ws.Range("F207").Value = 'Piast Gliwice'
$ws.Range("G207").Value = 'Radomiak Radom'
$ws.Range("K207").Value = 2
$ws.Range("L207").Value = 3.25
$ws.Range("M207").Value = 4
$ws.Range("N207").Value = 1.95
$ws.Range("O207").Value = 3.2
$ws.Range("P207").Value = 4.333
$ws.Range("Q207").Value = -0.5
$ws.Range("R207").Value = 1.95
$ws.Range("S207").Value = 1.9
$ws.Range("T207").Value = 2
$ws.Range("U207").Value = 1.9
$ws.Range("V207").Value = 1.95
$ws.Range("W207").Value = 0
$ws.Range("X207").Value = 0
$ws.Range("Y207").Value = 0
$ws.Range("Z207").Value = 0
$ws.Range("AA207").Value = 0

# Row 208
$ws.Range("A208").Value = 206
$ws.Range("B208").Value = 6775562
$ws.Range("C208").Value = 'Poland Ekstraklasa'
$ws.Range("D208").Value = 'Poland Ekstraklasa'
$ws.Range("E208").Value = 45359.6875
$ws.Range("F208").Value = 'Jagiellonia Bialystok'
$ws.Range("G208").Value = 'Slask Wroclaw'
$ws.Range("K208").Value = 1.952
$ws.Range("L208").Value = 3.5
$ws.Range("M208").Value = 3.8
$ws.Range("N208").Value = 1.909
$ws.Range("O208").Value = 3.5
$ws.Range("P208").Value = 3.8
$ws.Range("Q208").Value = -0.5
$ws.Range("R208").Value = 1.975
$ws.Range("S208").Value = 1.875
$ws.Range("T208").Value = 2.5
$ws.Range("U208").Value = 1.975
$ws.Range("V208").Value = 1.875
$ws.Range("W208").Value = 0
$ws.Range("X208").Value = 0
$ws.Range("Y208").Value = 0
$ws.Range("Z208").Value = 0
$ws.Range("AA208").Value = 0

# Row 209
$ws.Range("A209").Value = 207
$ws.Range("B209").Value = 6775563
$ws.Range("C209").Value = 'Poland Ekstraklasa'
$ws.Range("D209").Value = 'Poland Ekstraklasa'
$ws.Range("E209").Value = 45360.45833333334
$ws.Range("F209").Value = 'Korona Kielce'
$ws.Range("G209").Value = 'Cracovia Krakow'
$ws.Range("K209").Value = 2.625
$ws.Range("L209").Value = 3.2
$ws.Range("M209").Value = 2.7
$ws.Range("N209").Value = 2.4
$ws.Range("O209").Value = 3.2
$ws.Range("P209").Value = 2.9
$ws.Range("Q209").Value = -0.25
$ws.Range("R209").Value = 2.1
$ws.Range("S209").Value = 1.775
$ws.Range("T209").Value = 2.25
$ws.Range("U209").Value = 1.925
$ws.Range("V209").Value = 1.925
$ws.Range("W209").Value = 0
$ws.Range("X209").Value = 0
$ws.Range("Y209").Value = 0
$ws.Range("Z209").Value = 0
$ws.Range("AA209").Value = 0

# Row 210
$ws.Range("A210").Value = 208
$ws.Range("B210").Value = 6774879
$ws.Range("C210").Value = 'Poland Ekstraklasa'
$ws.Range("D210").Value = 'Poland Ekstraklasa'
$ws.Range("E210").Value = 45360.5625
$ws.Range("F210").Value = 'MKS Puszcza Niepolomice'
$ws.Range("G210").Value = 'Rakow Czestochowa'
$ws.Range("K210").Value = 4.75
$ws.Range("L210").Value = 3.8
$ws.Range("M210").Value = 1.7
$ws.Range("N210").Value = 5.25
$ws.Range("O210").Value = 4
$ws.Range("P210").Value = 1.615
$ws.Range("Q210").Value = 1
$ws.Range("R210").Value = 1.8
$ws.Range("S210").Value = 2.05
$ws.Range("T210").Value = 2.5
$ws.Range("U210").Value = 1.85
$ws.Range("V210").Value = 2
$ws.Range("W210").Value = 0
$ws.Range("X210").Value = 0
$ws.Range("Y210").Value = 0
$ws.Range("Z210").Value = 0
$ws.Range("AA210").Value = 0

# Row 211
$ws.Range("A211").Value = 209
$ws.Range("B211").Value = 6775561
$ws.Range("C211").Value = 'Poland Ekstraklasa'
$ws.Range("D211").Value = 'Poland Ekstraklasa'
$ws.Range("E211").Value = 45360.66666666666
$ws.Range("F211").Value = 'Gornik Zabrze'
$ws.Range("G211").Value = 'Lech Poznan'
$ws.Range("K211").Value = 3
$ws.Range("L211").Value = 3.4
$ws.Range("M211").Value = 2.3
$ws.Range("N211").Value = 2.9
$ws.Range("O211").Value = 3.4
$ws.Range("P211").Value = 2.375
$ws.Range("Q211").Value = 0.25
$ws.Range("R211").Value = 1.8
$ws.Range("S211").Value = 2.05
$ws.Range("T211").Value = 2.5
$ws.Range("U211").Value = 1.975
$ws.Range("V211").Value = 1.875
$ws.Range("W211").Value = 0
$ws.Range("X211").Value = 0
$ws.Range("Y211").Value = 0
$ws.Range("Z211").Value = 0
$ws.Range("AA211").Value = 0

# Row 212
$ws.Range("A212").Value = 210
$ws.Range("B212").Value = 6774464
$ws.Range("C212").Value = 'Poland Ekstraklasa'
$ws.Range("D212").Value = 'Poland Ekstraklasa'
$ws.Range("E212").Value = 45361.35416666666
$ws.Range("F212").Value = 'Stal Mielec'
$ws.Range("G212").Value = 'Ruch Chorzow'
$ws.Range("K212").Value = 2
$ws.Range("L212").Value = 3.4
$ws.Range("M212").Value = 3.75
$ws.Range("N212").Value = 2.55
$ws.Range("O212").Value = 3.1
$ws.Range("P212").Value = 2.8
$ws.Range("Q212").Value = 0
$ws.Range("R212").Value = 1.8
$ws.Range("S212").Value = 2.05
$ws.Range("T212").Value = 2.25
$ws.Range("U212").Value = 1.975
$ws.Range("V212").Value = 1.875
$ws.Range("W212").Value = 0
$ws.Range("X212").Value = 0
$ws.Range("Y212").Value = 0
$ws.Range("Z212").Value = 0
$ws.Range("AA212").Value = 0

# Row 213
$ws.Range("A213").Value = 211
$ws.Range("B213").Value = 6775565
$ws.Range("C213").Value = 'Poland Ekstraklasa'
$ws.Range("D213").Value = 'Poland Ekstraklasa'
$ws.Range("E213").Value = 45361.45833333334
$ws.Range("F213").Value = 'Pogon Szczecin'
$ws.Range("G213").Value = 'Zaglebie Lubin'
$ws.Range("K213").Value = 1.666
$ws.Range("L213").Value = 4
$ws.Range("M213").Value = 4.75
$ws.Range("N213").Value = 1.727
$ws.Range("O213").Value = 3.8
$ws.Range("P213").Value = 4.5
$ws.Range("Q213").Value = -0.75
$ws.Range("R213").Value = 1.975
$ws.Range("S213").Value = 1.875
$ws.Range("T213").Value = 2.75
$ws.Range("U213").Value = 1.9
$ws.Range("V213").Value = 1.95
$ws.Range("W213").Value = 0
$ws.Range("X213").Value = 0
$ws.Range("Y213").Value = 0
$ws.Range("Z213").Value = 0
$ws.Range("AA213").Value = 0

# Row 214
$ws.Range("A214").Value = 212
$ws.Range("B214").Value = 6775566
$ws.Range("C214").Value = 'Poland Ekstraklasa'
$ws.Range("D214").Value = 'Poland Ekstraklasa'
$ws.Range("E214").Value = 45361.5625
$ws.Range("F214").Value = 'Widzew Lodz'
$ws.Range("G214").Value = 'Legia Warsaw'
$ws.Range("K214").Value = 3.5
$ws.Range("L214").Value = 3.5
$ws.Range("M214").Value = 2.05
$ws.Range("N214").Value = 3.8
$ws.Range("O214").Value = 3.5
$ws.Range("P214").Value = 1.95
$ws.Range("Q214").Value = 0.5
$ws.Range("R214").Value = 1.85
$ws.Range("S214").Value = 2
$ws.Range("T214").Value = 2.5
$ws.Range("U214").Value = 1.925
$ws.Range("V214").Value = 1.925
$ws.Range("W214").Value = 0
$ws.Range("X214").Value = 0
$ws.Range("Y214").Value = 0
$ws.Range("Z214").Value = 0
$ws.Range("AA214").Value = 0

# Row 215
$ws.Range("A215").Value = 213
$ws.Range("B215").Value = 6774465
$ws.Range("C215").Value = 'Poland Ekstraklasa'
$ws.Range("D215").Value = 'Poland Ekstraklasa'
$ws.Range("E215").Value = 45362.625
$ws.Range("F215").Value = 'Warta Poznan'
$ws.Range("G215").Value = 'LKS Lodz'
$ws.Range("K215").Value = 1.95
$ws.Range("L215").Value = 3.3
$ws.Range("M215").Value = 4
$ws.Range("N215").Value = 1.95
$ws.Range("O215").Value = 3.3
$ws.Range("P215").Value = 4
$ws.Range("Q215").Value = -0.5
$ws.Range("R215").Value = 1.975
$ws.Range("S215").Value = 1.875
$ws.Range("T215").Value = 2.25
$ws.Range("U215").Value = 2
$ws.Range("V215").Value = 1.85
$ws.Range("W215").Value = 0
$ws.Range("X215").Value = 0
$ws.Range("Y215").Value = 0
$ws.Range("Z215").Value = 0
$ws.Range("AA215").Value = 0

